$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 5.828753
$ws.Range("H2").Value2 = 17.486259
$ws.Range("I2").Value2 = 0.1911291943607339
$ws.Range("J2").Value2 = 0.1911291943607339
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 0.3862596666666667
$ws.Range("N2").Value2 = 1.158779
$ws.Range("O2").Value2 = 0.05360826277999409
$ws.Range("P2").Value2 = 0.05360826277999409
$ws.Range("Q2").Value2 = 2.251412190862333
$ws.Range("R2").Value2 = 20.262709717761
$ws.Range("S2").Value2 = 0.01024610407621879
$ws.Range("T2").Value2 = 0.01024610407621879

$ws.Range("G3").Value2 = 5.828753
$ws.Range("H3").Value2 = 17.486259
$ws.Range("I3").Value2 = 0.1911291943607339
$ws.Range("J3").Value2 = 0.1911291943607339
$ws.Range("O3").Value2 = 0.1630272174193556
$ws.Range("P3").Value2 = 0.1630272174193557
$ws.Range("Q3").Value2 = 6.846733053944001
$ws.Range("R3").Value2 = 61.62059748549601
$ws.Range("S3").Value2 = 0.03115926072423365
$ws.Range("T3").Value2 = 0.03115926072423365

$ws.Range("G4").Value2 = 5.828753
$ws.Range("H4").Value2 = 17.486259
$ws.Range("I4").Value2 = 0.1911291943607339
$ws.Range("J4").Value2 = 0.1911291943607339
$ws.Range("M4").Value2 = 5.644318666666667
$ws.Range("N4").Value2 = 16.932956
$ws.Range("O4").Value2 = 0.7833645198006502
$ws.Range("P4").Value2 = 0.7833645198006502
$ws.Range("Q4").Value2 = 32.89933936128934
$ws.Range("R4").Value2 = 296.094054251604
$ws.Range("S4").Value2 = 0.1497238295602815
$ws.Range("T4").Value2 = 0.1497238295602815

$ws.Range("I5").Value2 = 0.7732994524709527
$ws.Range("J5").Value2 = 0.7732994524709526
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 0.3862596666666667
$ws.Range("N5").Value2 = 1.158779
$ws.Range("O5").Value2 = 0.05360826277999409
$ws.Range("P5").Value2 = 0.05360826277999409
$ws.Range("Q5").Value2 = 9.109104552569333
$ws.Range("R5").Value2 = 81.981940973124
$ws.Range("S5").Value2 = 0.04145524025568838
$ws.Range("T5").Value2 = 0.04145524025568838

$ws.Range("I6").Value2 = 0.7732994524709527
$ws.Range("J6").Value2 = 0.7732994524709526
$ws.Range("O6").Value2 = 0.1630272174193556
$ws.Range("P6").Value2 = 0.1630272174193557
$ws.Range("S6").Value2 = 0.1260688579682507
$ws.Range("T6").Value2 = 0.1260688579682507

$ws.Range("I7").Value2 = 0.7732994524709527
$ws.Range("J7").Value2 = 0.7732994524709526
$ws.Range("M7").Value2 = 5.644318666666667
$ws.Range("N7").Value2 = 16.932956
$ws.Range("O7").Value2 = 0.7833645198006502
$ws.Range("P7").Value2 = 0.7833645198006502
$ws.Range("Q7").Value2 = 133.1091317568373
$ws.Range("R7").Value2 = 1197.982185811536
$ws.Range("S7").Value2 = 0.6057753542470136
$ws.Range("T7").Value2 = 0.6057753542470135

$ws.Range("G8").Value2 = 1.084798333333333
$ws.Range("H8").Value2 = 3.254395
$ws.Range("I8").Value2 = 0.03557135316831352
$ws.Range("J8").Value2 = 0.03557135316831351
$ws.Range("K8").Value2 = 3
$ws.Range("L8").Value2 = 1
$ws.Range("M8").Value2 = 0.3862596666666667
$ws.Range("N8").Value2 = 1.158779
$ws.Range("O8").Value2 = 0.05360826277999409
$ws.Range("P8").Value2 = 0.05360826277999409
$ws.Range("Q8").Value2 = 0.4190138426338889
$ws.Range("R8").Value2 = 3.771124583705
$ws.Range("S8").Value2 = 0.001906918448086926
$ws.Range("T8").Value2 = 0.001906918448086926

$ws.Range("G9").Value2 = 1.084798333333333
$ws.Range("H9").Value2 = 3.254395
$ws.Range("I9").Value2 = 0.03557135316831352
$ws.Range("J9").Value2 = 0.03557135316831351
$ws.Range("O9").Value2 = 0.1630272174193556
$ws.Range("P9").Value2 = 0.1630272174193557
$ws.Range("Q9").Value2 = 1.274256192653334
$ws.Range("R9").Value2 = 11.46830573388
$ws.Range("S9").Value2 = 0.005799098726871332
$ws.Range("T9").Value2 = 0.005799098726871332

$ws.Range("G10").Value2 = 1.084798333333333
$ws.Range("H10").Value2 = 3.254395
$ws.Range("I10").Value2 = 0.03557135316831352
$ws.Range("J10").Value2 = 0.03557135316831351
$ws.Range("M10").Value2 = 5.644318666666667
$ws.Range("N10").Value2 = 16.932956
$ws.Range("O10").Value2 = 0.7833645198006502
$ws.Range("P10").Value2 = 0.7833645198006502
$ws.Range("Q10").Value2 = 6.122947482402222
$ws.Range("R10").Value2 = 55.10652734162001
$ws.Range("S10").Value2 = 0.02786533599335526
$ws.Range("T10").Value2 = 0.02786533599335525

